$d = $word.ActiveDocument

# ===========================================================================
# 1) Insert the four new "Ghost Scream" / "Jump Scare" paragraphs right after
#    the "To Be Recorded" paragraph (and before the following blank
#    paragraph that precedes "Rooms:").
# ===========================================================================

$toBeRecorded = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*To Be Recorded*") {
        $toBeRecorded = $p
    }
}

if ($toBeRecorded -eq $null) {
    throw "Could not locate the 'To Be Recorded' paragraph"
}

# Create a fresh, empty paragraph right after "To Be Recorded" by inserting a
# paragraph break before the (blank) paragraph that currently follows it.
$following = $toBeRecorded.Next()
$followingRange = $following.Range
$followingRange.Collapse(1)
[void]$followingRange.InsertParagraphBefore()

# That freshly created (still empty) paragraph now sits right after
# "To Be Recorded"; fill it (and push out three more paragraphs after it)
# via a single InsertXML call.
$newPara = $toBeRecorded.Next()

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:highlight w:val="green"/>
    </w:rPr>
    <w:t>Ghost Scream</w:t>
  </w:r>
  <w:r>
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Played w</w:t>
  </w:r>
  <w:r>
    <w:t>hen the ghost appears</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:firstLine="360"/>
    <w:rPr>
      <w:i/>
      <w:iCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t>Taken From Horror &amp; Suspense asset pack</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:highlight w:val="green"/>
    </w:rPr>
    <w:t>Jump Scare</w:t>
  </w:r>
  <w:r>
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Played w</w:t>
  </w:r>
  <w:r>
    <w:t>hen the ghost appears</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:firstLine="360"/>
    <w:rPr>
      <w:i/>
      <w:iCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t xml:space="preserve">Piano Jump Scare Stinger by </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t>TheSoundFXGuy_YT</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t xml:space="preserve"> -- https://freesound.org/s/534218/ -- License: Attribution 4.0</w:t>
  </w:r>
</w:p>
'@

[void]$newPara.Range.InsertXML($xml)

# ===========================================================================
# 2) Move the <w:lastRenderedPageBreak/> marker off of the "Kitchen" run and
#    onto the "Upstairs Bedroom" run. Both paragraphs are rewritten in full
#    (via InsertXML over their whole Range) with every attribute preserved,
#    since InsertXML replaces the complete <w:p> element of the range it is
#    applied to.
# ===========================================================================

$upstairsBedroom = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Upstairs Bedroom:*") {
        $upstairsBedroom = $p
    }
}

if ($upstairsBedroom -eq $null) {
    throw "Could not locate the 'Upstairs Bedroom' paragraph"
}

$upstairsBedroomRange = $upstairsBedroom.Range.Duplicate
$upstairsBedroomXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0325AC2D" w14:textId="4FF9852A" w:rsidR="006C4D29" w:rsidRDefault="006C4D29" w:rsidP="00A850A0">
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Upstairs Bedroom</w:t>
  </w:r>
  <w:r w:rsidRPr="008D7A4D">
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Player must ascend the stairs that go to the left. </w:t>
  </w:r>
  <w:r w:rsidR="00414A1A">
    <w:t xml:space="preserve">Enter the room immediately in front and the </w:t>
  </w:r>
  <w:r w:rsidR="00414A1A" w:rsidRPr="00827CDC">
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>cross</w:t>
  </w:r>
  <w:r w:rsidR="00414A1A">
    <w:t xml:space="preserve"> will be sitting on the bed.</w:t>
  </w:r>
  <w:r w:rsidR="004F1BBA">
    <w:t xml:space="preserve"> Light flickering sound plays on the ceiling lamp.</w:t>
  </w:r>
</w:p>
'@
[void]$upstairsBedroomRange.InsertXML($upstairsBedroomXml)

$kitchen = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Kitchen:*") {
        $kitchen = $p
    }
}

if ($kitchen -eq $null) {
    throw "Could not locate the 'Kitchen' paragraph"
}

$kitchenRange = $kitchen.Range.Duplicate
$kitchenXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0323FAB0" w14:textId="6C5B06F0" w:rsidR="00B07093" w:rsidRDefault="00B07093" w:rsidP="00A850A0">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Kitchen</w:t>
  </w:r>
  <w:r w:rsidRPr="008D7A4D">
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Tangentially connected to the main entrance + grand staircase. </w:t>
  </w:r>
  <w:r w:rsidRPr="00827CDC">
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Holy water</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> can be found in this room</w:t>
  </w:r>
  <w:r w:rsidR="00827CDC">
    <w:t xml:space="preserve"> on the kitchen </w:t>
  </w:r>
  <w:r w:rsidR="00C24233">
    <w:t>countertop</w:t>
  </w:r>
  <w:r>
    <w:t>.</w:t>
  </w:r>
  <w:r w:rsidR="002512C3">
    <w:t xml:space="preserve"> Lots of metal objects will cause the sound to reflect more than in other rooms.</w:t>
  </w:r>
</w:p>
'@
[void]$kitchenRange.InsertXML($kitchenXml)

Write-Host "Edit applied: inserted Ghost Scream / Jump Scare entries and relocated lastRenderedPageBreak to Upstairs Bedroom."
